# Form101.xlsx smoke-test update:
#  - fix three header labels in row 1 (AG1/AJ1/AK1) that were showing
#    sample answer text instead of the field name
#  - fix a few stray cell values in row 2 (AE2/AG2/AJ2/AN2)
#  - duplicate row 2 into a brand-new row 3 (a second QA test-case row),
#    tweaking a handful of the copied fields
#  - update the sheet view (selection / scroll) and used-range dimension

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- row 1: header label corrections -------------------------------------
$ws.Range("AG1").Value = "CauseofInjury"
$ws.Range("AJ1").Value = "Doing"
$ws.Range("AK1").Value = "DOLastHire"

# ---- row 2: a few cell corrections ----------------------------------------
$ws.Range("AE2").Value = "Head - Eyes"
$ws.Range("AG2").Value = "Motor Vehicle - Vehicle Upset"
$ws.Range("AJ2").Value = "TestDoing"
$ws.Range("AN2").Value = "test@gmail.com"

# ---- row 3: duplicate row 2, then adjust a handful of fields --------------
# Copy formats first (so number formats land on the existing date style
# instead of minting a new one), then copy values on top.
$ws.Range("A2:AN2").Copy()
$ws.Range("A3:AN3").PasteSpecial(-4122)
$ws.Range("A2:AN2").Copy()
$ws.Range("A3:AN3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("A3").Value = "TC002"
$ws.Range("B3").Value = "Test "
$ws.Range("F3").Value = "Test "
$ws.Range("G3").Value = "Test "
$ws.Range("J3").Value = "Married"
$ws.Range("M3").Value = "QA"
$ws.Range("AN3").Value = "test"

# ---- sheet view: selection + used range ------------------------------------
$ws.Range("AL9").Select()
$excel.ActiveWindow.ScrollColumn = 25
$excel.ActiveWindow.ScrollRow = 1
